$d = $word.ActiveDocument

# Update the date heading.
$d.Content.Find.Execute("2024-03-26 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-03-27 Wednesday", 2) | Out-Null

# Update the multiplication problems table, cell by cell (row, column) so
# that the duplicate "79x12=" entries are disambiguated correctly.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "76×76="
$t.Cell(1, 2).Range.Text = "74×48="
$t.Cell(1, 3).Range.Text = "79×12="
$t.Cell(1, 4).Range.Text = "36×39="
$t.Cell(1, 5).Range.Text = "21×44="

$t.Cell(5, 1).Range.Text = "41×81="
$t.Cell(5, 2).Range.Text = "39×33="
$t.Cell(5, 3).Range.Text = "20×22="
$t.Cell(5, 4).Range.Text = "33×73="
$t.Cell(5, 5).Range.Text = "46×23="

$t.Cell(10, 1).Range.Text = "91×32="
$t.Cell(10, 2).Range.Text = "51×35="
$t.Cell(10, 3).Range.Text = "53×76="
$t.Cell(10, 4).Range.Text = "25×33="
$t.Cell(10, 5).Range.Text = "39×36="

$t.Cell(15, 1).Range.Text = "61×13="
$t.Cell(15, 2).Range.Text = "27×65="
$t.Cell(15, 3).Range.Text = "67×12="
$t.Cell(15, 4).Range.Text = "18×97="
$t.Cell(15, 5).Range.Text = "87×55="

$t.Cell(20, 1).Range.Text = "94×91="
$t.Cell(20, 2).Range.Text = "34×98="
$t.Cell(20, 3).Range.Text = "27×93="
$t.Cell(20, 4).Range.Text = "17×27="
$t.Cell(20, 5).Range.Text = "83×42="
